$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3966.794699136231
$ws.Range("C2").Value = 10.55041620778134
$ws.Range("D2").Value = 2210.467911692193
